$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices / percentages) remain as text
# by pre-setting the cell NumberFormat to "@" (Text) before assigning .Value
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

# Apply updated cell values row by row
# Row 2
$ws.Range("D2").Value = '247.19'
$ws.Range("E2").Value = '0.81%'
# Row 3
$ws.Range("D3").Value = '26.26'
$ws.Range("E3").Value = '4.18%'
# Row 4
$ws.Range("D4").Value = '5.087'
$ws.Range("E4").Value = '1.40%'
# Row 5
$ws.Range("D5").Value = '0.05600'
$ws.Range("E5").Value = '-0.23%'
# Row 6
$ws.Range("D6").Value = '6.486'
$ws.Range("E6").Value = '-1.26%'
# Row 7
$ws.Range("D7").Value = '0.8129'
$ws.Range("E7").Value = '-0.05%'
# Row 8
$ws.Range("D8").Value = '0.8440'
$ws.Range("E8").Value = '0.76%'
# Row 9
$ws.Range("B9").Value = 'MandalaExchangeToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D9").Value = '0.07004'
$ws.Range("E9").Value = '0.82%'
# Row 10
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D10").Value = '0.02811'
$ws.Range("E10").Value = '-1.16%'
# Row 11
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D11").Value = '0.09397'
$ws.Range("E11").Value = '-0.04%'
# Row 12
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D12").Value = '0.001516'
$ws.Range("E12").Value = '-0.32%'
# Row 13
$ws.Range("B13").Value = 'One'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D13").Value = '0.0006009'
$ws.Range("E13").Value = '0.91%'
# Row 14
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").Value = '0.006192'
$ws.Range("E14").Value = '0.42%'
# Row 15
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").Value = '3.607'
$ws.Range("E15").Value = '3.08%'
# Row 16
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D16").Value = '3.020'
$ws.Range("E16").Value = '0.44%'
# Row 17
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").Value = '2.055'
$ws.Range("E17").Value = '-1.75%'
# Row 18
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").Value = '0.3126'
$ws.Range("E18").Value = '-1.81%'
# Row 19
$ws.Range("B19").Value = 'WazirX'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D19").Value = '0.1341'
$ws.Range("E19").Value = '0.29%'
# Row 20
$ws.Range("D20").Value = '0.03179'
$ws.Range("E20").Value = '-2.61%'
# Row 21
$ws.Range("E21").Value = '-1.30%'
# Row 22
$ws.Range("D22").Value = '3.751'
$ws.Range("E22").Value = '0.40%'
# Row 23
$ws.Range("D23").Value = '0.04658'
$ws.Range("E23").Value = '-0.07%'
# Row 25
$ws.Range("D25").Value = '0.001243'
$ws.Range("E25").Value = '-0.03%'
# Row 26
$ws.Range("D26").Value = '0.004579'
$ws.Range("E26").Value = '17.56%'
# Row 27
$ws.Range("D27").Value = '0.00009599'
$ws.Range("E27").Value = '-0.93%'
# Row 28
$ws.Range("E28").Value = '-0.05%'
# Row 40
$ws.Range("D40").Value = '0.03665'
$ws.Range("E40").Value = '0.05%'
# Row 41
$ws.Range("D41").Value = '0.006149'
$ws.Range("E41").Value = '81.40%'
# Row 42
$ws.Range("D42").Value = '0.1057'
$ws.Range("E42").Value = '-22.49%'
# Row 43
$ws.Range("D43").Value = '0.002500'
$ws.Range("E43").Value = '-8.32%'
# Row 44
$ws.Range("D44").Value = '0.008242'
$ws.Range("E44").Value = '0.58%'
# Row 45
$ws.Range("D45").Value = '0.00005450'
$ws.Range("E45").Value = '2.97%'
# Row 46
$ws.Range("E46").Value = '0.06%'
# Row 47
$ws.Range("D47").Value = '0.1450'
$ws.Range("E47").Value = '-35.80%'
# Row 48
$ws.Range("D48").Value = '0.002586'
$ws.Range("E48").Value = '26.87%'
# Row 49
$ws.Range("E49").Value = '0.06%'
# Row 50
$ws.Range("E50").Value = '0.06%'
